# "ff kijke of het kan" - move the source table from rows 13-17 up to
# rows 1-5, re-point the chart's series/category/value references at the
# new location, restyle the chart and reposition it on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Shift the little results table from A13:C17 up to A1:C5 -----------
# Deleting the now-unused rows above the table shifts everything (values,
# formulas, shared-string refs) up without touching the data itself.
$ws.Rows("1:12").Delete() | Out-Null

# --- 2. Re-point the chart series formulas at the table's new location ----
$co = $ws.ChartObjects(1)
$chart = $co.Chart

$s1 = $chart.SeriesCollection(1)
$s1.Formula = "=SERIES(Blad1!`$A`$2,Blad1!`$B`$1:`$C`$1,Blad1!`$B`$2:`$C`$2,1)"

$s2 = $chart.SeriesCollection(2)
$s2.Formula = "=SERIES(Blad1!`$A`$3,Blad1!`$B`$1:`$C`$1,Blad1!`$B`$3:`$C`$3,2)"

$s3 = $chart.SeriesCollection(3)
$s3.Formula = "=SERIES(Blad1!`$A`$4,Blad1!`$B`$1:`$C`$1,Blad1!`$B`$4:`$C`$4,3)"

$s4 = $chart.SeriesCollection(4)
$s4.Formula = "=SERIES(Blad1!`$A`$5,Blad1!`$B`$1:`$C`$1,Blad1!`$B`$5:`$C`$5,4)"

# --- 3. Apply the new chart style (Design > Chart Styles > Style 7) -------
$chart.ChartStyle = 7

# --- 4. Move/resize the chart to its new anchor ----------------------------
$co.Left = 67.06238188976378
$co.Top = 116.99984251968505
$co.Width = 466.37503937007875
$co.Height = 225.00007874015748

# --- 5. Reset the view: scroll back to A1 and select G4 --------------------
$ws.Activate() | Out-Null
$ws.Range("A1").Select() | Out-Null
$ws.Range("G4").Select() | Out-Null
